$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "university" building block (rows 97-101), matching the pattern of
# the other building sections (townhall, builder, residence, hospital, ...).
# Building name goes in column A, start/end positions in C/D; column E is
# the pre-existing shared "structurize scan" formula that recalculates once
# C/D/A/B are populated.

# Row 97
$ws.Range("A97").Value = "university"
$ws.Range("C97").Value = "'1480 -51 1"
$ws.Range("D97").Value = "'1508 -30 -36"

# Row 98
$ws.Range("A98").Value = "university"
$ws.Range("C98").Value = "1480 -51 -47"
$ws.Range("D98").Value = "1508 -30 -84"

# Row 99
$ws.Range("A99").Value = "university"
$ws.Range("C99").Value = "'1480 -51 -95"
$ws.Range("D99").Value = "1508 -30 -132"

# Row 100
$ws.Range("A100").Value = "university"
$ws.Range("C100").Value = "'1480 -51 -143"
$ws.Range("D100").Value = "'1508 -30 -180"

# Row 101
$ws.Range("A101").Value = "university"
$ws.Range("C101").Value = "'1480 -51 -191"
$ws.Range("D101").Value = "'1508 -30 -228"

# Match the author's final selection/scroll position in the sheet.
$ws.Range("A85").Select()
$ws.Range("E97").Select()
